$wb = $excel.ActiveWorkbook

function Set-TextLabel {
    param($range)
    # Excel auto-converts a purely-numeric string (e.g. "2050") back into
    # a number when assigned through .Value, even though the target cell
    # must stay a text label. Force text entry via Text number format,
    # then restore the original (General/bordered-header) formatting by
    # copying it back from a neighboring cell that already holds the
    # correct style, so only the *value*'s type changes, not its style.
    $fmtSource = $range.Offset(0, -1)
    $range.NumberFormat = "@"
    $range.Value = "2050"
    $fmtSource.Copy()
    $range.PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false
}

# Sheets whose header row E1 label needs fixing from the stray numeric
# value to the correct "2050" text label, and whose trailing "Total" row
# (row 13) must be removed.
$sheetsYearLabel = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)"
)

foreach ($name in $sheetsYearLabel) {
    $ws = $wb.Worksheets.Item($name)
    Set-TextLabel $ws.Range("E1")
    $ws.Rows.Item(13).Delete()
}

# "Potencia Incremental - SIN(MW)" uses period ranges instead of single
# years, so its E1 label is "2041-2050" instead of "2050". That text is
# not purely numeric, so a plain assignment already stores it as text.
$ws = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$ws.Range("E1").Value = "2041-2050"
$ws.Rows.Item(13).Delete()

# "Emissoes Totais (MtCO2eq)" only needs the E1 label fix - it never had
# a Total row.
$ws = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")
Set-TextLabel $ws.Range("E1")

# "Custo Total (bilhões de R$)" has no year-label row at all, just its
# own trailing "Total" row (row 4) to remove.
$ws = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws.Rows.Item(4).Delete()

$wb.Save()
